$wb = $excel.ActiveWorkbook

# --- Create the new "classes" sheet by duplicating an existing sheet so it
#     inherits the exact same style palette (fonts/fills/borders/alignment) ---
$src = $wb.Worksheets.Item("classrooms")
$src.Copy($wb.Worksheets.Item(1))
$ws = $wb.Worksheets.Item(1)
$ws.Name = "classes"

# Remove merges inherited from the copied sheet; they will be recreated below
$ws.Cells.UnMerge()

# The template sheet had 34 data rows; "classes" only needs 24
$ws.Range("A25:F34").EntireRow.Delete()

# "classes" uses a wider "names" column than the copied template
$ws.Columns.Item(5).ColumnWidth = 13

# --- Re-stripe the B:F formatting per group (odd groups / even groups) ---
$ws.Range("B3:F4").Copy()
$ws.Range("B10:F17").PasteSpecial(-4122)
$ws.Range("B5:F9").Copy()
$ws.Range("B18:F19").PasteSpecial(-4122)
$ws.Range("B3:F4").Copy()
$ws.Range("B20:F24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Header row ---
$ws.Cells.Item(2,2).Value = 'group_No.'
$ws.Cells.Item(2,3).Value = 'names_base'
$ws.Cells.Item(2,4).Value = 'names_in_group_No.'
$ws.Cells.Item(2,5).Value = 'names'
$ws.Cells.Item(2,6).Value = 'names_No.'

# Group 1
$ws.Cells.Item(3,2).Value = '1.'
$ws.Cells.Item(3,3).Value = 1
$ws.Cells.Item(3,4).Value = '1.'
$ws.Cells.Item(3,5).Value = '1B 1gas_log'
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(4,4).Value = '2.'
$ws.Cells.Item(4,5).Value = '1d 1wz'
$ws.Cells.Item(4,6).Value = 2
$ws.Range("B3:B4").Merge()
$ws.Range("C3:C4").Merge()

# Group 2
$ws.Cells.Item(5,2).Value = '2.'
$ws.Cells.Item(5,3).Value = 2
$ws.Cells.Item(5,4).Value = '1.'
$ws.Cells.Item(5,5).Value = '2A 2ek_fry'
$ws.Cells.Item(5,6).Value = 3
$ws.Cells.Item(6,4).Value = '2.'
$ws.Cells.Item(6,5).Value = '2B 2gas'
$ws.Cells.Item(6,6).Value = 4
$ws.Cells.Item(7,4).Value = '3.'
$ws.Cells.Item(7,5).Value = '2C 2hot'
$ws.Cells.Item(7,6).Value = 5
$ws.Cells.Item(8,4).Value = '4.'
$ws.Cells.Item(8,5).Value = '2L 2log'
$ws.Cells.Item(8,6).Value = 6
$ws.Cells.Item(9,4).Value = '5.'
$ws.Cells.Item(9,5).Value = '2d 2wz'
$ws.Cells.Item(9,6).Value = 7
$ws.Range("B5:B9").Merge()
$ws.Range("C5:C9").Merge()

# Group 3
$ws.Cells.Item(10,2).Value = '3.'
$ws.Cells.Item(10,3).Value = 3
$ws.Cells.Item(10,4).Value = '1.'
$ws.Cells.Item(10,5).Value = '3A 3ra_fry'
$ws.Cells.Item(10,6).Value = 8
$ws.Cells.Item(11,4).Value = '2.'
$ws.Cells.Item(11,5).Value = '3B 3gas'
$ws.Cells.Item(11,6).Value = 9
$ws.Cells.Item(12,4).Value = '3.'
$ws.Cells.Item(12,5).Value = '3C 3hot'
$ws.Cells.Item(12,6).Value = 10
$ws.Cells.Item(13,4).Value = '4.'
$ws.Cells.Item(13,5).Value = '3JL'
$ws.Cells.Item(13,6).Value = 11
$ws.Cells.Item(14,4).Value = '5.'
$ws.Cells.Item(14,5).Value = '3L 3log'
$ws.Cells.Item(14,6).Value = 12
$ws.Cells.Item(15,4).Value = '6.'
$ws.Cells.Item(15,5).Value = '3S 3sport'
$ws.Cells.Item(15,6).Value = 13
$ws.Cells.Item(16,4).Value = '7.'
$ws.Cells.Item(16,5).Value = '3d 3wz'
$ws.Cells.Item(16,6).Value = 14
$ws.Cells.Item(17,4).Value = '8.'
$ws.Cells.Item(17,5).Value = '3e 3wz'
$ws.Cells.Item(17,6).Value = 15
$ws.Range("B10:B17").Merge()
$ws.Range("C10:C17").Merge()

# Group 4
$ws.Cells.Item(18,2).Value = '4.'
$ws.Cells.Item(18,3).Value = 4
$ws.Cells.Item(18,4).Value = '1.'
$ws.Cells.Item(18,5).Value = '4A 4ra_log'
$ws.Cells.Item(18,6).Value = 16
$ws.Cells.Item(19,4).Value = '2.'
$ws.Cells.Item(19,5).Value = '4B 4gas_fry'
$ws.Cells.Item(19,6).Value = 17
$ws.Range("B18:B19").Merge()
$ws.Range("C18:C19").Merge()

# Group 5
$ws.Cells.Item(20,2).Value = '5.'
$ws.Cells.Item(20,3).Value = 5
$ws.Cells.Item(20,4).Value = '1.'
$ws.Cells.Item(20,5).Value = '5A 5ek_log'
$ws.Cells.Item(20,6).Value = 18
$ws.Cells.Item(21,4).Value = '2.'
$ws.Cells.Item(21,5).Value = '5B 5gas_fry'
$ws.Cells.Item(21,6).Value = 19
$ws.Cells.Item(22,4).Value = '3.'
$ws.Cells.Item(22,5).Value = '5MP'
$ws.Cells.Item(22,6).Value = 20
$ws.Cells.Item(23,4).Value = '4.'
$ws.Cells.Item(23,5).Value = '5OP'
$ws.Cells.Item(23,6).Value = 21
$ws.Cells.Item(24,4).Value = '5.'
$ws.Cells.Item(24,5).Value = '5ZI'
$ws.Cells.Item(24,6).Value = 22
$ws.Range("B20:B24").Merge()
$ws.Range("C20:C24").Merge()

$ws.Range("A1").Select()
